# Updated cryptos list on Sun Sep 10 02:52:26 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 45/46 swapped places (Cronos <-> RenderToken); handle B/C/D/E fully for those rows.
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.49"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.61%  "

$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0520"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.22%  "

# Price (D) and Volume(1h) (E) updates for the remaining rows.
$ws.Range("D2").Value = "25.976.02"
$ws.Range("E2").Value = "  +0.15%  "

$ws.Range("D3").Value = "1.638.73"
$ws.Range("E3").Value = "  +0.01%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.01"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.33%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.18%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.506"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.17%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.01"
$ws.Range("D7").Style = "Normal"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.0635"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.80%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.253"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.40%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.40"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.21%  "

$ws.Range("E11").Value = "  -0.15%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.24"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.34%  "

$ws.Range("D13").Value = "1.650.60"
$ws.Range("E13").Value = "  +0.64%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.540"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.79%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "63.08"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.20%  "

$ws.Range("D16").Value = "0.0₃0757"
$ws.Range("E16").Value = "  -1.22%  "

$ws.Range("D17").Value = "25.991.89"
$ws.Range("E17").Value = "  +0.18%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.01"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.32%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "193.30"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.12%  "

$ws.Range("E20").Value = "  -1.33%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.72"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.11%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.15"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.30%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.132"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.25%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "143.98"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.08%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.01"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.27%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.77"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.06%  "

$ws.Range("E27").Value = "  -0.57%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.45"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.64%  "

$ws.Range("E29").Value = "  -0.32%  "

$ws.Range("E30").Value = "  -3.08%  "

$ws.Range("E31").Value = "  +0.34%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.26"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.60%  "

$ws.Range("E33").Value = "  -1.17%  "

$ws.Range("E34").Value = "  +0.45%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.896"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.68%  "

$ws.Range("D36").Value = "1.122.79"
$ws.Range("E36").Value = "  -1.53%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.46"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.19%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.533"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.34%  "

$ws.Range("E39").Value = "  -1.31%  "

$ws.Range("E40").Value = "  -0.18%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "98.23"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.18%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.30"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.22%  "

$ws.Range("E43").Value = "  -0.21%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "56.03"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.05%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.69"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.32%  "

$ws.Range("E48").Value = "  +0.04%  "

$ws.Range("E49").Value = "  +0.14%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0939"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.01%  "

$ws.Range("E51").Value = "  +0.49%  "

